$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2030.25
$ws.Range("I33").Value = 891.7143
$ws.Range("K33").Value = 891.7143
$ws.Range("M33").Value = -662.7143
$ws.Range("H55").Value = 63391.375
$ws.Range("J55").Value = 72401.57000000001
$ws.Range("L55").Value = 72401.57000000001
$ws.Range("N55").Value = -72829.57000000001
$ws.Range("H80").Value = 809.5
$ws.Range("I80").Value = 1107.2858
$ws.Range("J80").Value = 577.8889
$ws.Range("K80").Value = 3321.8574
$ws.Range("L80").Value = 1733.6667
$ws.Range("M80").Value = -2323.8574
$ws.Range("N80").Value = -3729.6667
$ws.Range("H83").Value = 809.5
$ws.Range("I83").Value = 1107.2858
$ws.Range("J83").Value = 577.8889
$ws.Range("K83").Value = 9965.572200000001
$ws.Range("L83").Value = 5201.0001
$ws.Range("M83").Value = -4973.572200000001
$ws.Range("N83").Value = -15185.0001
$ws.Range("H86").Value = 2925896.5
$ws.Range("I86").Value = 2223.25
$ws.Range("K86").Value = 2223.25
$ws.Range("M86").Value = -1100.25
$ws.Range("H89").Value = 2925896.5
$ws.Range("I89").Value = 2223.25
$ws.Range("K89").Value = 11116.25
$ws.Range("M89").Value = -5500.25
$ws.Range("H106").Value = 3472.625
$ws.Range("I106").Value = 4457.6924
$ws.Range("J106").Value = 2308.4546
$ws.Range("K106").Value = 4457.6924
$ws.Range("L106").Value = 2308.4546
$ws.Range("M106").Value = -3826.6924
$ws.Range("N106").Value = -3570.4546
$ws.Range("H111").Value = 40246.617
$ws.Range("I111").Value = 60366.41
$ws.Range("J111").Value = 2242.5557
$ws.Range("K111").Value = 181099.23
$ws.Range("L111").Value = 6727.6671
$ws.Range("M111").Value = -178032.23
$ws.Range("N111").Value = -12861.6671
$ws.Range("H118").Value = 63511.668
$ws.Range("I118").Value = 63511.668
$ws.Range("K118").Value = 190535.004
$ws.Range("M118").Value = -188878.004
$ws.Range("H128").Value = 89315.375
$ws.Range("J128").Value = 89315.375
$ws.Range("L128").Value = 89315.375
$ws.Range("N128").Value = -99275.375
$ws.Range("H129").Value = 16579.8
$ws.Range("J129").Value = 12288.363
$ws.Range("L129").Value = 36865.089
$ws.Range("N129").Value = -46865.089
$ws.Range("H132").Value = 1472180.9
$ws.Range("J132").Value = 4385.273
$ws.Range("L132").Value = 13155.819
$ws.Range("N132").Value = -18215.819
$ws.Range("H133").Value = 59613.77
$ws.Range("J133").Value = 59613.77
$ws.Range("L133").Value = 59613.77
$ws.Range("N133").Value = -69733.76999999999
$ws.Range("H137").Value = 3696.6428
$ws.Range("J137").Value = 6149.4116
$ws.Range("L137").Value = 18448.2348
$ws.Range("N137").Value = -23548.2348
$ws.Range("H138").Value = 3982.2856
$ws.Range("I138").Value = 2140.5366
$ws.Range("J138").Value = 5738.372
$ws.Range("K138").Value = 6421.6098
$ws.Range("L138").Value = 17215.116
$ws.Range("M138").Value = -1281.6098
$ws.Range("N138").Value = -27495.116

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("I42").Value = 20000
$ws.Range("K42").Value = 20000
$ws.Range("M42").Value = -19514
$ws.Range("H74").Value = 1990.1212
$ws.Range("I74").Value = 1868.069
$ws.Range("K74").Value = 1868.069
$ws.Range("M74").Value = -994.069
$ws.Range("H77").Value = 1990.1212
$ws.Range("I77").Value = 1868.069
$ws.Range("K77").Value = 9340.344999999999
$ws.Range("M77").Value = -4972.344999999999
$ws.Range("H97").Value = 1703.24
$ws.Range("I97").Value = 1997.3158
$ws.Range("K97").Value = 1997.3158
$ws.Range("M97").Value = -1501.3158
$ws.Range("H132").Value = 4540.0757
$ws.Range("I132").Value = 1876.5714
$ws.Range("K132").Value = 5629.7142
$ws.Range("M132").Value = -3099.7142

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3036.524
$ws.Range("I107").Value = 2856.4285
$ws.Range("J107").Value = 3396.7144
$ws.Range("K107").Value = 2856.4285
$ws.Range("L107").Value = 3396.7144
$ws.Range("M107").Value = -936.4285
$ws.Range("N107").Value = -7236.7144
$ws.Range("H134").Value = 17476.203
$ws.Range("I134").Value = 2019.0754
$ws.Range("J134").Value = 68677.94
$ws.Range("K134").Value = 6057.2262
$ws.Range("L134").Value = 206033.82
$ws.Range("M134").Value = -3522.2262
$ws.Range("N134").Value = -211103.82

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 265.8889
$ws.Range("J22").Value = 247
$ws.Range("L22").Value = 247
$ws.Range("N22").Value = -947
$ws.Range("H31").Value = 225414.62
$ws.Range("I31").Value = 371815.97
$ws.Range("K31").Value = 371815.97
$ws.Range("M31").Value = -371520.97
$ws.Range("H34").Value = 225414.62
$ws.Range("I34").Value = 371815.97
$ws.Range("K34").Value = 371815.97
$ws.Range("M34").Value = -371613.97
$ws.Range("H86").Value = 8616.333000000001
$ws.Range("I86").Value = 7509.6
$ws.Range("K86").Value = 7509.6
$ws.Range("M86").Value = -6386.6
$ws.Range("H89").Value = 8616.333000000001
$ws.Range("I89").Value = 7509.6
$ws.Range("K89").Value = 37548
$ws.Range("M89").Value = -31932
$ws.Range("H131").Value = 64649.5
$ws.Range("I131").Value = 49999
$ws.Range("J131").Value = 79300
$ws.Range("K131").Value = 49999
$ws.Range("L131").Value = 79300
$ws.Range("M131").Value = -44959
$ws.Range("N131").Value = -89380
$ws.Range("H132").Value = 3475.6956
$ws.Range("I132").Value = 2654.5806
$ws.Range("K132").Value = 7963.7418
$ws.Range("M132").Value = -5433.7418
$ws.Range("H134").Value = 195992.33
$ws.Range("I134").Value = 2509.139
$ws.Range("K134").Value = 7527.417
$ws.Range("M134").Value = -4992.417

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 3087130
$ws.Range("J113").Value = 732.25
$ws.Range("L113").Value = 2196.75
$ws.Range("N113").Value = -6536.75
$ws.Range("H118").Value = 2569.1924
$ws.Range("I118").Value = 449.75
$ws.Range("J118").Value = 2954.5454
$ws.Range("K118").Value = 1349.25
$ws.Range("L118").Value = 8863.636200000001
$ws.Range("M118").Value = -106.25
$ws.Range("N118").Value = -11349.6362

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H18").Value = 333391680
$ws.Range("J18").Value = 333391680
$ws.Range("L18").Value = 333391680
$ws.Range("N18").Value = -333392266
$ws.Range("H63").Value = 20114
$ws.Range("J63").Value = 20114
$ws.Range("L63").Value = 20114
$ws.Range("N63").Value = -21486
$ws.Range("H66").Value = 20114
$ws.Range("J66").Value = 20114
$ws.Range("L66").Value = 60342
$ws.Range("N66").Value = -67206
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("H97").Value = 4263.032
$ws.Range("I97").Value = 5373.7393
$ws.Range("K97").Value = 5373.7393
$ws.Range("M97").Value = -4877.7393
$ws.Range("H134").Value = 49998.332
$ws.Range("J134").Value = 49998.332
$ws.Range("L134").Value = 149994.996
$ws.Range("N134").Value = -155064.996
$ws.Range("H141").Value = 52500
$ws.Range("J141").Value = 52500
$ws.Range("L141").Value = 52500
$ws.Range("N141").Value = -62860
$ws.Range("N69").ClearContents()
$ws.Range("N72").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 20001
$ws.Range("I2").Value = 20001
$ws.Range("K2").Value = 20001
$ws.Range("M2").Value = -19889
$ws.Range("H20").Value = 100
$ws.Range("J20").Value = 100
$ws.Range("L20").Value = 100
$ws.Range("N20").Value = -552
$ws.Range("H42").Value = 14964
$ws.Range("J42").Value = 14964
$ws.Range("L42").Value = 14964
$ws.Range("N42").Value = -16090
$ws.Range("H46").Value = 4864.095
$ws.Range("I46").Value = 4354.1665
$ws.Range("K46").Value = 4354.1665
$ws.Range("M46").Value = -4166.1665
$ws.Range("H49").Value = 14964
$ws.Range("J49").Value = 14964
$ws.Range("L49").Value = 14964
$ws.Range("N49").Value = -15258
$ws.Range("H128").Value = 75000
$ws.Range("J128").Value = 75000
$ws.Range("L128").Value = 75000
$ws.Range("N128").Value = -84960
$ws.Range("H132").Value = 5142.472
$ws.Range("I132").Value = 4891.9375
$ws.Range("J132").Value = 5342.9
$ws.Range("K132").Value = 14675.8125
$ws.Range("L132").Value = 16028.7
$ws.Range("M132").Value = -12145.8125
$ws.Range("N132").Value = -21088.7
$ws.Range("H141").Value = 52500
$ws.Range("J141").Value = 52500
$ws.Range("L141").Value = 52500
$ws.Range("N141").Value = -62860

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 67240.39999999999
$ws.Range("I107").Value = 83804.664
$ws.Range("K107").Value = 251413.992
$ws.Range("M107").Value = -249493.992
$ws.Range("H126").Value = 4588
$ws.Range("I126").Value = 4333.8
$ws.Range("K126").Value = 13001.4
$ws.Range("M126").Value = -10531.4
$ws.Range("H132").Value = 25857.627
$ws.Range("I132").Value = 1577.1072
$ws.Range("K132").Value = 4731.321599999999
$ws.Range("M132").Value = -2201.321599999999
$ws.Range("H136").Value = 32919.27
$ws.Range("I136").Value = 923.5179000000001
$ws.Range("K136").Value = 2770.5537
$ws.Range("M136").Value = -220.5537000000004
